$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text, matching the original inline-string cells,
# so Excel does not auto-convert numeric-looking values (e.g. "248.69") into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.034.41"
$ws.Range("E2").Value = "  +4.66%  "
$ws.Range("D3").Value = "1.916.81"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "248.69"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "0.686"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "47.61"
$ws.Range("E8").Value = "  +10.06%  "
$ws.Range("D9").Value = "0.375"
$ws.Range("E9").Value = "  +5.66%  "
$ws.Range("D10").Value = "58.21"
$ws.Range("E10").Value = "  +6.20%  "
$ws.Range("D11").Value = "0.0758"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "15.52"
$ws.Range("E13").Value = "  +12.52%  "
$ws.Range("D14").Value = "0.823"
$ws.Range("E14").Value = "  +6.80%  "
$ws.Range("D15").Value = "2.192.39"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "5.12"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "1.920.78"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "36.992.19"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").Value = "74.55"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "13.63"
$ws.Range("E21").Value = "  +6.28%  "
$ws.Range("D22").Value = "250.40"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").Value = "5.14"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("D26").Value = "167.36"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").Value = "8.79"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").Value = "18.65"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "4.54"
$ws.Range("E31").Value = "  +5.98%  "
$ws.Range("D32").Value = "0.0608"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("E33").Value = "  +26.41%  "
$ws.Range("D34").Value = "4.28"
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("D35").Value = "1.89"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "18.98"
$ws.Range("E37").Value = "  +37.76%  "
$ws.Range("D38").Value = "0.888"
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("D39").Value = "1.45"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D42").Value = "0.0226"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").Value = "17.60"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("E44").Value = "  +19.70%  "
$ws.Range("D45").Value = "1.09"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "1.349.74"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("D48").Value = "0.0834"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("D49").Value = "2.80"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("D50").Value = "6.38"
$ws.Range("E50").Value = "  +1.73%  "

# Row 40/41 swap: Aave <-> LidoDAOToken (coin ranking changed order)
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "105.90"
$ws.Range("E41").Value = "  +8.37%  "

# Row 51: RocketPoolETH replaced with THORChain
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "3.75"
$ws.Range("E51").Value = "  +13.08%  "

# Restore default styling on the Price column (remove the temporary text format)
$priceRange.Style = "Normal"
